$d = $word.ActiveDocument
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷2=34, 0", 2)
$d.Content.Find.Execute("41÷4=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2)
$d.Content.Find.Execute("94÷7=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2)
$d.Content.Find.Execute("54÷5=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "97÷6=16, 1", 2)
$d.Content.Find.Execute("42÷4=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "95÷4=23, 3", 2)
$d.Content.Find.Execute("64÷4=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=17, 0", 2)
$d.Content.Find.Execute("45÷9=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "63÷5=12, 3", 2)
$d.Content.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷3=8, 1", 2)
$d.Content.Find.Execute("69÷7=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=11, 7", 2)
$d.Content.Find.Execute("79÷7=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "62÷5=12, 2", 2)
$d.Content.Find.Execute("97÷7=13, 6", $true, $false, $false, $false, $false, $true, 1, $false, "84÷5=16, 4", 2)
$d.Content.Find.Execute("72÷9=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "98÷2=49, 0", 2)
$d.Content.Find.Execute("26÷2=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2)
$d.Content.Find.Execute("94÷8=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "84÷8=10, 4", 2)
$d.Content.Find.Execute("96÷8=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=9, 7", 2)
$d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=6, 6", 2)
$d.Content.Find.Execute("57÷9=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=8, 4", 2)
$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "74÷6=12, 2", 2)
$d.Content.Find.Execute("49÷6=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "85÷3=28, 1", 2)
$d.Content.Find.Execute("45÷7=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "87÷2=43, 1", 2)
$d.Content.Find.Execute("30÷7=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=1, 2", 2)
$d.Content.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "65÷9=7, 2", 2)
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "42÷3=14, 0", 2)
$d.Content.Find.Execute("19÷5=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "69÷5=13, 4", 2)
$d.Content.Find.Execute("28÷7=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷7=10, 5", 2)
